$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2023-01-05, serial 44931) is inserted above the
# existing row 398, shifting the remaining "Pimiento" rows (old 398-422) down
# by one (new 399-423).
$ws.Rows.Item(398).Insert()

$ws.Cells.Item(398, 1).Value = 11
$ws.Cells.Item(398, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(398, 3).Value = "Bíobío"
$ws.Cells.Item(398, 4).Value = 44931
$ws.Cells.Item(398, 5).Value = 8
$ws.Cells.Item(398, 6).Value = 100112002
$ws.Cells.Item(398, 7).Value = "Pimiento"
$ws.Cells.Item(398, 8).Value = "Zafiro rojo"
$ws.Cells.Item(398, 9).Value = "Primera"
$ws.Cells.Item(398, 10).Value = 180
$ws.Cells.Item(398, 11).Value = 17000
$ws.Cells.Item(398, 12).Value = 18000
$ws.Cells.Item(398, 13).Value = 17556
$ws.Cells.Item(398, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(398, 15).Value = "Limache"
$ws.Cells.Item(398, 16).Value = 975
$ws.Cells.Item(398, 17).Value = 18
$ws.Cells.Item(398, 18).Value = "Hortaliza"
